$d = $word.ActiveDocument

# --- Rename the use-case list items (paragraphs 7-13 in the original order) ---

# 7: Kreiranje racuna -> Kreiranje dokumenata
$p = $d.Paragraphs.Item(7)
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "Kreiranje dokumenata"

# 8: Kreiranje izdatnice -> Pregled dokumenata
$p = $d.Paragraphs.Item(8)
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "Pregled dokumenata"

# 9: Kreiranje primke -> "Upravlj" + bookmark(_GoBack) + "anje resursima"
$p = $d.Paragraphs.Item(9)
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "Upravljanje resursima"
$findRng = $d.Content
[void]$findRng.Find.Execute("Upravlj", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRng = $d.Range($findRng.End, $findRng.End)
[void]$d.Bookmarks.Add("_GoBack", $bmRng)

# 10: Kreiranje narudzbenice -> Pregled resursa
$p = $d.Paragraphs.Item(10)
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "Pregled resursa"

# 11: Pretrazivanje zaposlenika -> Slanje maila
$p = $d.Paragraphs.Item(11)
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "Slanje maila"

# 12: Pretrazivanje opreme -> Kreiranje izvjestaja
$p = $d.Paragraphs.Item(12)
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "Kreiranje izvještaja"

# 13: Pretrazivanje artikala -> Obavijest o manjku zaliha
$p = $d.Paragraphs.Item(13)
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "Obavijest o manjku zaliha"

# --- Remove the now-duplicated trailing paragraphs (old items 14, 15, 16) ---
# Delete from the highest index down so earlier indices stay valid.

$p = $d.Paragraphs.Item(16)
$p.Range.Delete()

$p = $d.Paragraphs.Item(15)
$p.Range.Delete()

$p = $d.Paragraphs.Item(14)
$p.Range.Delete()

Write-Output "done"
